$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 143295.2310000001
$ws.Range("C2").Value = 754978.4700000002
$ws.Range("D2").Value = 32550.93899999999
$ws.Range("E2").Value = 599667.1729999994
$ws.Range("F2").Value = 187179.1249999998
$ws.Range("G2").Value = 119355.1630000001
$ws.Range("B3").Value = 165679.6790000002
$ws.Range("C3").Value = 711866.3900000006
$ws.Range("D3").Value = 15753.226
$ws.Range("E3").Value = 628411.0559999974
$ws.Range("F3").Value = 185809.68
$ws.Range("G3").Value = 110499.72
$ws.Range("H3").Value = 167395.916
$ws.Range("C4").Value = 773660.0140000001
$ws.Range("D4").Value = 46460.832
$ws.Range("E4").Value = 638320.2009999987
$ws.Range("G4").Value = 117631.867
$ws.Range("H4").Value = 148023.1900000001
$ws.Range("B5").Value = 170960.6570000002
$ws.Range("C5").Value = 724488.0519999998
$ws.Range("D5").Value = 60139.187
$ws.Range("E5").Value = 718431.0040000005
$ws.Range("F5").Value = 301521.8129999997
$ws.Range("G5").Value = 185381.8519999998
$ws.Range("B6").Value = 181767.2709999999
$ws.Range("C6").Value = 844929.1660000004
$ws.Range("D6").Value = 57862.82800000001
$ws.Range("E6").Value = 704934.6430000024
$ws.Range("F6").Value = 331113.1789999998
$ws.Range("G6").Value = 192459.0310000004
$ws.Range("H6").Value = 142973.212
$ws.Range("C7").Value = 780127.389
$ws.Range("D7").Value = 81598.87699999999
$ws.Range("E7").Value = 815701.4070000015
$ws.Range("F7").Value = 295500.6100000001
$ws.Range("G7").Value = 402969.7210000003
$ws.Range("H7").Value = 159794.0620000001
